$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 481.15
$ws.Range("I15").Value = 481.15
$ws.Range("K15").Value = 1443.45
$ws.Range("M15").Value = -1274.45
# Row 17
$ws.Range("H17").Value = 1452.1063
$ws.Range("J17").Value = 1452.1063
$ws.Range("L17").Value = 4356.3189
$ws.Range("N17").Value = -4692.3189
# Row 64
$ws.Range("H64").Value = 3676.8533
$ws.Range("J64").Value = 3993.125
$ws.Range("L64").Value = 3993.125
$ws.Range("N64").Value = -4489.125
# Row 67
$ws.Range("H67").Value = 3676.8533
$ws.Range("J67").Value = 3993.125
$ws.Range("L67").Value = 3993.125
$ws.Range("N67").Value = -5709.125
# Row 76
$ws.Range("H76").Value = 3140.878
$ws.Range("I76").Value = 2709.1428
$ws.Range("K76").Value = 2709.1428
$ws.Range("M76").Value = -2394.1428
# Row 79
$ws.Range("H79").Value = 3140.878
$ws.Range("I79").Value = 2709.1428
$ws.Range("K79").Value = 2709.1428
$ws.Range("M79").Value = -1617.1428
# Row 112
$ws.Range("H112").Value = 2411.2222
$ws.Range("J112").Value = 2650.125
$ws.Range("L112").Value = 7950.375
$ws.Range("N112").Value = -10166.375
# Row 137
$ws.Range("H137").Value = 7829684.5
$ws.Range("I137").Value = 12195923
$ws.Range("J137").Value = 46389.652
$ws.Range("K137").Value = 36587769
$ws.Range("L137").Value = 139168.956
$ws.Range("M137").Value = -36585219
$ws.Range("N137").Value = -144268.956
# Row 138
$ws.Range("H138").Value = 2436.2805
$ws.Range("I138").Value = 2472.5
$ws.Range("J138").Value = 2432.3647
$ws.Range("K138").Value = 7417.5
$ws.Range("L138").Value = 7297.0941
$ws.Range("M138").Value = -2277.5
$ws.Range("N138").Value = -17577.0941

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2247.805
$ws.Range("I61").Value = 2149.0303
$ws.Range("J61").Value = 2655.25
$ws.Range("K61").Value = 2149.0303
$ws.Range("L61").Value = 2655.25
$ws.Range("M61").Value = -1937.0303
$ws.Range("N61").Value = -3079.25
# Row 74
$ws.Range("H74").Value = 1301.62
$ws.Range("I74").Value = 1077.5435
$ws.Range("K74").Value = 1077.5435
$ws.Range("M74").Value = -203.5435
# Row 76
$ws.Range("H76").Value = 20192
$ws.Range("J76").Value = 20192
$ws.Range("L76").Value = 20192
$ws.Range("N76").Value = -20868
# Row 77
$ws.Range("H77").Value = 1301.62
$ws.Range("I77").Value = 1077.5435
$ws.Range("K77").Value = 5387.7175
$ws.Range("M77").Value = -1019.7175
# Row 79
$ws.Range("H79").Value = 20192
$ws.Range("J79").Value = 20192
$ws.Range("L79").Value = 20192
$ws.Range("N79").Value = -22532
# Row 132
$ws.Range("H132").Value = 1776.1852
$ws.Range("I132").Value = 909.1875
$ws.Range("J132").Value = 3037.2727
$ws.Range("K132").Value = 2727.5625
$ws.Range("L132").Value = 9111.8181
$ws.Range("M132").Value = -197.5625
$ws.Range("N132").Value = -14171.8181
# Row 136
$ws.Range("H136").Value = 2247.805
$ws.Range("I136").Value = 2149.0303
$ws.Range("J136").Value = 2655.25
$ws.Range("K136").Value = 6447.090899999999
$ws.Range("L136").Value = 7965.75
$ws.Range("M136").Value = -3897.090899999999
$ws.Range("N136").Value = -13065.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1356.6154
$ws.Range("I134").Value = 935.30304
$ws.Range("J134").Value = 3673.8333
$ws.Range("K134").Value = 2805.90912
$ws.Range("L134").Value = 11021.4999
$ws.Range("M134").Value = -270.9091200000003
$ws.Range("N134").Value = -16091.4999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2423.2593
$ws.Range("I58").Value = 1845.7333
$ws.Range("J58").Value = 3145.1667
$ws.Range("K58").Value = 1845.7333
$ws.Range("L58").Value = 3145.1667
$ws.Range("M58").Value = -1642.7333
$ws.Range("N58").Value = -3551.1667
# Row 99
$ws.Range("H99").Value = 2944.7778
$ws.Range("I99").Value = 2166.6667
$ws.Range("J99").Value = 3333.8333
$ws.Range("K99").Value = 2166.6667
$ws.Range("L99").Value = 3333.8333
$ws.Range("M99").Value = -668.6667000000002
$ws.Range("N99").Value = -6329.8333
# Row 107
$ws.Range("H107").Value = 670.5
$ws.Range("I107").Value = 390.33334
$ws.Range("J107").Value = 838.6
$ws.Range("K107").Value = 390.33334
$ws.Range("L107").Value = 838.6
$ws.Range("M107").Value = 1529.66666
$ws.Range("N107").Value = -4678.6
# Row 126
$ws.Range("H126").Value = 2944.7778
$ws.Range("I126").Value = 2166.6667
$ws.Range("J126").Value = 3333.8333
$ws.Range("K126").Value = 6500.000100000001
$ws.Range("L126").Value = 10001.4999
$ws.Range("M126").Value = -4030.000100000001
$ws.Range("N126").Value = -14941.4999
# Row 132
$ws.Range("H132").Value = 2220.2856
$ws.Range("I132").Value = 2007.1154
$ws.Range("J132").Value = 2836.111
$ws.Range("K132").Value = 6021.3462
$ws.Range("L132").Value = 8508.332999999999
$ws.Range("M132").Value = -3491.3462
$ws.Range("N132").Value = -13568.333
# Row 134
$ws.Range("H134").Value = 2216
$ws.Range("I134").Value = 1264.6666
$ws.Range("J134").Value = 3396.9656
$ws.Range("K134").Value = 3793.9998
$ws.Range("L134").Value = 10190.8968
$ws.Range("M134").Value = -1258.9998
$ws.Range("N134").Value = -15260.8968
# Row 136
$ws.Range("H136").Value = 2423.2593
$ws.Range("I136").Value = 1845.7333
$ws.Range("J136").Value = 3145.1667
$ws.Range("K136").Value = 5537.199900000001
$ws.Range("L136").Value = 9435.500100000001
$ws.Range("M136").Value = -2987.199900000001
$ws.Range("N136").Value = -14535.5001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 3458.25
$ws.Range("I59").Value = 2500
$ws.Range("J59").Value = 3777.6667
$ws.Range("K59").Value = 7500
$ws.Range("L59").Value = 11333.0001
$ws.Range("M59").Value = -6960
$ws.Range("N59").Value = -12413.0001
# Row 60
$ws.Range("H60").Value = 302.875
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 302.875
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 908.625
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -1410.625
# Row 76
$ws.Range("H76").Value = 3966.6667
$ws.Range("I76").Value = 1933.3334
$ws.Range("J76").Value = 4983.3335
$ws.Range("K76").Value = 5800.0002
$ws.Range("L76").Value = 14950.0005
$ws.Range("M76").Value = -5417.0002
$ws.Range("N76").Value = -15716.0005
# Row 79
$ws.Range("H79").Value = 3966.6667
$ws.Range("I79").Value = 1933.3334
$ws.Range("J79").Value = 4983.3335
$ws.Range("K79").Value = 5800.0002
$ws.Range("L79").Value = 14950.0005
$ws.Range("M79").Value = -4474.0002
$ws.Range("N79").Value = -17602.0005
# Row 131
$ws.Range("H131").Value = 883.10254
$ws.Range("I131").Value = 680
$ws.Range("J131").Value = 894.08105
$ws.Range("K131").Value = 2040
$ws.Range("L131").Value = 2682.24315
$ws.Range("M131").Value = 3000
$ws.Range("N131").Value = -12762.24315

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 1547.7778
$ws.Range("I126").Value = 1987.2
$ws.Range("J126").Value = 998.5
$ws.Range("K126").Value = 5961.6
$ws.Range("L126").Value = 2995.5
$ws.Range("M126").Value = -3491.6
$ws.Range("N126").Value = -7935.5
# Row 132
$ws.Range("H132").Value = 3867.25
$ws.Range("I132").Value = 4105.852
$ws.Range("J132").Value = 3371.6924
$ws.Range("K132").Value = 12317.556
$ws.Range("L132").Value = 10115.0772
$ws.Range("M132").Value = -9787.556
$ws.Range("N132").Value = -15175.0772

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 136
$ws.Range("H136").Value = 2257.5
$ws.Range("I136").Value = 1593.75
$ws.Range("J136").Value = 4381.5
$ws.Range("K136").Value = 4781.25
$ws.Range("L136").Value = 13144.5
$ws.Range("M136").Value = -2231.25
$ws.Range("N136").Value = -18244.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2510.1904
$ws.Range("I122").Value = 2448.7693
$ws.Range("J122").Value = 2610
$ws.Range("K122").Value = 7346.3079
$ws.Range("L122").Value = 7830
$ws.Range("M122").Value = -4896.3079
$ws.Range("N122").Value = -12730
# Row 132
$ws.Range("H132").Value = 2017.8049
$ws.Range("I132").Value = 1077.9565
$ws.Range("J132").Value = 3218.7222
$ws.Range("K132").Value = 3233.8695
$ws.Range("L132").Value = 9656.1666
$ws.Range("M132").Value = -703.8694999999998
$ws.Range("N132").Value = -14716.1666
